{"js": "// The title paragraph currently reads:\n//   \"PODIUM DES COUPES DU MONDE ET QUELQUES CONTEXTES ECONOMIQUES ET POLITIQUES\"\n// It must become:\n//   \"EVOLUTION D\\u2019UN PAYS DANS LE CARRE FINAL LIEE AUX CONTEXTES ECONOMIQUES ET POLITIQUES\"\n//\n// i.e. \"PODIUM DES COUPES DU MONDE\" -> \"EVOLUTION D\\u2019UN PAYS DANS LE CARRE FINAL LIEE AUX\"\n// and  \" ET QUELQUES \"             -> \" \" (single space)\n\nconst body = context.document.body;\n\n// Replace the leading phrase, preserving the run's (bold, size 32) formatting.\nconst titleMatches = body.search(\"PODIUM DES COUPES DU MONDE\", { matchCase: true, matchWholeWord: false });\ntitleMatches.load(\"items\");\nawait context.sync();\n\nif (titleMatches.items.length === 0) {\n  throw new Error(\"Could not find the title text 'PODIUM DES COUPES DU MONDE'.\");\n}\ntitleMatches.items[0].insertText(\"EVOLUTION D\\u2019UN PAYS DANS LE CARRE FINAL LIEE AUX\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Collapse \" ET QUELQUES \" down to a single space, removing the \"QUELQUES \" run.\nconst midMatches = body.search(\" ET QUELQUES \", { matchCase: true, matchWholeWord: false });\nmidMatches.load(\"items\");\nawait context.sync();\n\nif (midMatches.items.length === 0) {\n  throw new Error(\"Could not find the ' ET QUELQUES ' text to collapse.\");\n}\nmidMatches.items[0].insertText(\" \", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The title paragraph currently reads:\n#   \"PODIUM DES COUPES DU MONDE ET QUELQUES CONTEXTES ECONOMIQUES ET POLITIQUES\"\n# It must become:\n#   \"EVOLUTION D'UN PAYS DANS LE CARRE FINAL LIEE AUX CONTEXTES ECONOMIQUES ET POLITIQUES\"\n# (using a right single quotation mark, U+2019, in \"D'UN\")\n#\n# i.e. \"PODIUM DES COUPES DU MONDE ET QUELQUES CONTEXTES\" becomes\n#      \"EVOLUTION D\\u2019UN PAYS DANS LE CARRE FINAL LIEE AUX CONTEXTES\"\n# which both swaps the leading phrase and collapses \" ET QUELQUES \" down to \" \".\n\n$d = $word.ActiveDocument\n\n$apostrophe = [char]0x2019\n$oldPhrase = \"PODIUM DES COUPES DU MONDE ET QUELQUES CONTEXTES\"\n$newPhrase = \"EVOLUTION D\" + $apostrophe + \"UN PAYS DANS LE CARRE FINAL LIEE AUX CONTEXTES\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldPhrase\n$find.Replacement.Text = $newPhrase\n\n$found = $find.Execute([ref]$find.Text, [ref]$true, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$find.Replacement.Text, [ref]2)\n\nif (-not $found) {\n    throw \"Could not find the title text to replace.\"\n}\n"}
